$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.482.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.447.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.50%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.78%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.150"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.23%  "
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("E11").Value = "  -4.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.900.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "68.408.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("E15").Value = "  -4.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.442.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.91%  "
$ws.Range("E18").Value = "  -2.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "342.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.33%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.00%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.571.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0837"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.34%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  -3.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "431.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +103.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -5.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("E41").Value = "  -4.33%  "
$ws.Range("E42").Value = "  -4.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("E45").Value = "  -6.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.32%  "
$ws.Range("E47").Value = "  -3.26%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0718"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.488"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.561"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0913"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.70%  "
